# Apply the authored changes to the Input_S&T / Input_AC sheets:
#  - update several "Selected_Value" numbers on Input_S&T (sheet 1)
#  - change the WF Position value from "pipe" to "shell"
#  - append a new "Tube update" / "update nt" row (row 34)
#  - move the active sheet / selection / scroll position from
#    Input_AC (sheet 2) back to Input_S&T (sheet 1)

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Input_S&T")
$ws2 = $wb.Worksheets.Item("Input_AC")

# --- Input_S&T (sheet1) cell value updates ---------------------------------

$ws1.Range("E2").Value = "shell"

$ws1.Range("E3").Value = 0.3
$ws1.Range("E4").Value = 101325
$ws1.Range("E5").Value = 363.15
$ws1.Range("E6").Value = 343.15

$ws1.Range("E11").Value = 303.15
$ws1.Range("E12").Value = 311.15

$ws1.Range("E16").Value = 785
$ws1.Range("E17").Value = 50000
$ws1.Range("E18").Value = 50000

$ws1.Range("E24").Value = 2
$ws1.Range("E25").Value = 0.5

$ws1.Range("E28").Value = 1.35
$ws1.Range("E29").Value = 0.5
$ws1.Range("E30").Value = 0.333
$ws1.Range("E31").Value = 30

# --- New row 34: "Tube update" / "update nt" -------------------------------
# Copy formatting from existing same-styled cells, then set the values so the
# new cells line up with the rest of the table (B -> centered/bordered label
# style, E -> centered value style).

[void]$ws1.Range("B4").Copy()
[void]$ws1.Range("B34").PasteSpecial(-4122)
$ws1.Range("B34").Value = "Tube update"

[void]$ws1.Range("D33").Copy()
[void]$ws1.Range("E34").PasteSpecial(-4122)
$ws1.Range("E34").Value = "update nt"

# --- View / selection: active tab moves from Input_AC to Input_S&T --------

[void]$ws2.Activate()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws2.Range("H15").Select()

[void]$ws1.Activate()
$excel.ActiveWindow.ScrollRow = 16
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws1.Range("F23").Select()
